# Regenerate "想去人数" (want-to-go count) / "最低票价" (min ticket price)
# figures, as published by the gh-pages data-refresh job (commit 456a3b4).
#
# Sheet "展览" (Exhibitions)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value  = 426          # F2  424 -> 426
$ws.Cells.Item(2, 7).Value  = "已售罄"      # G2   58 -> sold out
$ws.Cells.Item(3, 6).Value  = 416          # F3  415 -> 416
$ws.Cells.Item(4, 6).Value  = 2704         # F4  2698 -> 2704
$ws.Cells.Item(5, 6).Value  = 1319         # F5  1318 -> 1319
$ws.Cells.Item(7, 6).Value  = 1967         # F7  1968 -> 1967
$ws.Cells.Item(9, 6).Value  = 39           # F9  38 -> 39
$ws.Cells.Item(10, 6).Value = 580          # F10 579 -> 580
$ws.Cells.Item(11, 6).Value = 272          # F11 271 -> 272
$ws.Cells.Item(13, 6).Value = 11233        # F13 11215 -> 11233
$ws.Cells.Item(14, 6).Value = 6411         # F14 6399 -> 6411
$ws.Cells.Item(18, 6).Value = 249          # F18 247 -> 249
$ws.Cells.Item(21, 6).Value = 886          # F21 883 -> 886
$ws.Cells.Item(22, 6).Value = 36           # F22 35 -> 36
$ws.Cells.Item(23, 6).Value = 238          # F23 237 -> 238
$ws.Cells.Item(25, 6).Value = 3614         # F25 3612 -> 3614
$ws.Cells.Item(33, 6).Value = 4957         # F33 4955 -> 4957
$ws.Cells.Item(35, 6).Value = 1210         # F35 1208 -> 1210
$ws.Cells.Item(36, 6).Value = 199          # F36 198 -> 199
$ws.Cells.Item(37, 6).Value = 382          # F37 377 -> 382
$ws.Cells.Item(38, 6).Value = 165          # F38 163 -> 165

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(11, 6).Value = 119          # F11 118 -> 119
$ws.Cells.Item(12, 6).Value = 3643         # F12 3641 -> 3643

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 8948          # F2 8947 -> 8948
$ws.Cells.Item(3, 6).Value = 474           # F3 472 -> 474
$ws.Cells.Item(4, 6).Value = 1763          # F4 1759 -> 1763

# Sheet "全部类型" (All types - union of the above)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value  = 8948         # F2  8947 -> 8948
$ws.Cells.Item(3, 6).Value  = 474          # F3  472 -> 474
$ws.Cells.Item(4, 6).Value  = 1763         # F4  1759 -> 1763
$ws.Cells.Item(6, 6).Value  = 416          # F6  415 -> 416
$ws.Cells.Item(7, 6).Value  = 2704         # F7  2698 -> 2704
$ws.Cells.Item(10, 6).Value = 1319         # F10 1318 -> 1319
$ws.Cells.Item(14, 6).Value = 39           # F14 38 -> 39
$ws.Cells.Item(16, 6).Value = 580          # F16 579 -> 580
$ws.Cells.Item(17, 6).Value = 272          # F17 271 -> 272
$ws.Cells.Item(19, 6).Value = 11233        # F19 11215 -> 11233
$ws.Cells.Item(20, 6).Value = 3643         # F20 3641 -> 3643
$ws.Cells.Item(21, 6).Value = 6411         # F21 6399 -> 6411
$ws.Cells.Item(26, 6).Value = 249          # F26 247 -> 249
$ws.Cells.Item(29, 6).Value = 886          # F29 883 -> 886
$ws.Cells.Item(30, 6).Value = 36           # F30 35 -> 36
$ws.Cells.Item(31, 6).Value = 238          # F31 237 -> 238
$ws.Cells.Item(33, 6).Value = 3614         # F33 3612 -> 3614
$ws.Cells.Item(41, 6).Value = 4957         # F41 4955 -> 4957
$ws.Cells.Item(43, 6).Value = 1210         # F43 1208 -> 1210
$ws.Cells.Item(45, 6).Value = 165          # F45 163 -> 165
